$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1982921009669723
$ws.Range("B3").Value = 0.2118068965517241
$ws.Range("B4").Value = 0.1805458229957766
$ws.Range("B8").Value = 0.2757242757242757
$ws.Range("B12").Value = 0.216893039049236
$ws.Range("B13").Value = 0.2738805263656158
$ws.Range("B16").Value = 0.2832591683289857
